$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.992.07'
$ws.Range("E2").Value = '  +4.45%  '
$ws.Range("D3").Value = '2.701.67'
$ws.Range("E3").Value = '  +3.78%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.60'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.71'
$ws.Range("E6").Value = '  +3.87%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.606'
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("D9").Value = '2.729.87'
$ws.Range("E9").Value = '  +4.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.71'
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("E11").Value = '  +7.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.387'
$ws.Range("E12").Value = '  +4.46%  '
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("D14").Value = '3.182.93'
$ws.Range("E14").Value = '  +3.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.46'
$ws.Range("E15").Value = '  +8.14%  '
$ws.Range("D16").Value = '62.897.33'
$ws.Range("E16").Value = '  +4.26%  '
$ws.Range("E17").Value = '  +7.15%  '
$ws.Range("D18").Value = '2.721.10'
$ws.Range("E18").Value = '  +4.23%  '
$ws.Range("E19").Value = '  +5.20%  '
$ws.Range("E20").Value = '  +5.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '362.58'
$ws.Range("E21").Value = '  +4.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.98'
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.38'
$ws.Range("E25").Value = '  +2.76%  '
$ws.Range("E26").Value = '  +3.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.59'
$ws.Range("E27").Value = '  +7.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.02'
$ws.Range("E29").Value = '  +6.18%  '
$ws.Range("D30").Value = '0.0₃0848'
$ws.Range("E30").Value = '  +6.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.08'
$ws.Range("E31").Value = '  +10.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '169.42'
$ws.Range("E32").Value = '  +1.49%  '
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.53'
$ws.Range("E34").Value = '  +5.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.18'
$ws.Range("E35").Value = '  +19.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.74'
$ws.Range("E37").Value = '  +7.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("E38").Value = '  +10.17%  '
$ws.Range("E39").Value = '  +20.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '351.75'
$ws.Range("E40").Value = '  +12.96%  '
$ws.Range("E41").Value = '  +9.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.06'
$ws.Range("E42").Value = '  +2.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.63'
$ws.Range("E43").Value = '  +13.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.60'
$ws.Range("E44").Value = '  +8.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0596'
$ws.Range("E45").Value = '  +8.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.59'
$ws.Range("E46").Value = '  +8.63%  '
$ws.Range("E47").Value = '  +7.17%  '
$ws.Range("E48").Value = '  +5.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.18'
$ws.Range("E49").Value = '  +0.84%  '
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("E51").Value = '  -0.38%  '
